$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "68.153.54"
Set-TextValue $ws.Range("E2") "  +1.49%  "

Set-TextValue $ws.Range("D3") "3.910.47"
Set-TextValue $ws.Range("E3") "  -0.42%  "

Set-TextValue $ws.Range("E4") "  +0.10%  "

Set-TextValue $ws.Range("D5") "490.42"
Set-TextValue $ws.Range("E5") "  +4.18%  "

Set-TextValue $ws.Range("D6") "146.78"
Set-TextValue $ws.Range("E6") "  +0.57%  "

Set-TextValue $ws.Range("D7") "0.622"
Set-TextValue $ws.Range("E7") "  -0.84%  "

Set-TextValue $ws.Range("D9") "0.730"
Set-TextValue $ws.Range("E9") "  -0.52%  "

Set-TextValue $ws.Range("D10") "0.165"
Set-TextValue $ws.Range("E10") "  -0.13%  "

Set-TextValue $ws.Range("D11") "0.0000344"
Set-TextValue $ws.Range("E11") "  +1.14%  "

Set-TextValue $ws.Range("D12") "42.91"
Set-TextValue $ws.Range("E12") "  -1.22%  "

Set-TextValue $ws.Range("D13") "10.80"
Set-TextValue $ws.Range("E13") "  +3.39%  "

Set-TextValue $ws.Range("D14") "4.535.37"
Set-TextValue $ws.Range("E14") "  -0.48%  "

Set-TextValue $ws.Range("D15") "3.922.13"
Set-TextValue $ws.Range("E15") "  -1.32%  "

Set-TextValue $ws.Range("D16") "14.14"
Set-TextValue $ws.Range("E16") "  -6.85%  "

Set-TextValue $ws.Range("E17") "  -1.26%  "

Set-TextValue $ws.Range("D18") "19.84"
Set-TextValue $ws.Range("E18") "  -0.04%  "

Set-TextValue $ws.Range("E19") "  -2.09%  "

Set-TextValue $ws.Range("D20") "68.327.78"
Set-TextValue $ws.Range("E20") "  +1.30%  "

Set-TextValue $ws.Range("D21") "438.49"
Set-TextValue $ws.Range("E21") "  +0.05%  "

Set-TextValue $ws.Range("B22") "ImmutableX"
Set-TextValue $ws.Range("C22") "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D22") "3.53"
Set-TextValue $ws.Range("E22") "  +4.37%  "

Set-TextValue $ws.Range("B23") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C23") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D23") "15.04"
Set-TextValue $ws.Range("E23") "  +3.33%  "

Set-TextValue $ws.Range("D24") "87.78"
Set-TextValue $ws.Range("E24") "  -0.11%  "

Set-TextValue $ws.Range("D25") "11.49"
Set-TextValue $ws.Range("E25") "  +18.39%  "

Set-TextValue $ws.Range("D26") "11.43"
Set-TextValue $ws.Range("E26") "  +10.74%  "

Set-TextValue $ws.Range("D27") "3.65"
Set-TextValue $ws.Range("E27") "  +0.71%  "

Set-TextValue $ws.Range("D28") "38.32"
Set-TextValue $ws.Range("E28") "  -1.28%  "

Set-TextValue $ws.Range("D29") "5.75"
Set-TextValue $ws.Range("E29") "  -0.30%  "

Set-TextValue $ws.Range("D30") "726.34"
Set-TextValue $ws.Range("E30") "  +0.76%  "

Set-TextValue $ws.Range("D31") "13.70"
Set-TextValue $ws.Range("E31") "  +1.03%  "

Set-TextValue $ws.Range("E32") "  -1.28%  "

Set-TextValue $ws.Range("E33") "  +3.02%  "

Set-TextValue $ws.Range("E34") "  +17.40%  "

Set-TextValue $ws.Range("D35") "41.74"
Set-TextValue $ws.Range("E35") "  -2.64%  "

Set-TextValue $ws.Range("D36") "60.37"
Set-TextValue $ws.Range("E36") "  +4.20%  "

Set-TextValue $ws.Range("D37") "0.0₃0856"
Set-TextValue $ws.Range("E37") "  +7.30%  "

Set-TextValue $ws.Range("D38") "0.412"
Set-TextValue $ws.Range("E38") "  +22.39%  "

Set-TextValue $ws.Range("D39") "0.148"
Set-TextValue $ws.Range("E39") "  -1.85%  "

Set-TextValue $ws.Range("D40") "0.998"
Set-TextValue $ws.Range("E40") "  -0.11%  "

Set-TextValue $ws.Range("D41") "2.96"
Set-TextValue $ws.Range("E41") "  +15.16%  "

Set-TextValue $ws.Range("B42") "VeChain"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D42") "0.0480"
Set-TextValue $ws.Range("E42") "  +0.59%  "

Set-TextValue $ws.Range("B43") "ThetaToken"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D43") "3.15"
Set-TextValue $ws.Range("E43") "  +3.03%  "

Set-TextValue $ws.Range("E44") "  +3.84%  "

Set-TextValue $ws.Range("E45") "  -0.28%  "

Set-TextValue $ws.Range("E46") "  +0.12%  "

Set-TextValue $ws.Range("D47") "3.27"
Set-TextValue $ws.Range("E47") "  +3.50%  "

Set-TextValue $ws.Range("D48") "3.41"
Set-TextValue $ws.Range("E48") "  -3.97%  "

Set-TextValue $ws.Range("B49") "ARBITRUM"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D49") "2.13"
Set-TextValue $ws.Range("E49") "  -3.35%  "

Set-TextValue $ws.Range("B50") "Monero"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D50") "145.26"
Set-TextValue $ws.Range("E50") "  -1.45%  "

Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₆0337"
Set-TextValue $ws.Range("E51") "  +27.14%  "
